# Weekly update: a new price record (row) is reported for this market/product
# combination. It is inserted as the new row 13, pushing the previously
# existing rows 13-29 down to rows 14-30 (dimension grows from R29 to R30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts rows 13..29 down to 14..30).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record's data.
$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C13").Value = "Coquimbo"
$ws.Range("D13").Value = 44483
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 100112032
$ws.Range("G13").Value = "Zapallo italiano"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 340
$ws.Range("K13").Value = 10000
$ws.Range("L13").Value = 11000
$ws.Range("M13").Value = 10500
$ws.Range("N13").Value = "$/caja 60 unidades"
$ws.Range("O13").Value = "Provincia de Limarí"
$ws.Range("P13").Value = 175
$ws.Range("Q13").Value = 60
$ws.Range("R13").Value = "Hortaliza"
